$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"0.1361090093851089"
$ws.Cells.Item(2, 2).Value = [double]"0.9589046835899353"
$ws.Cells.Item(2, 3).Value = [double]"0.01270137634128332"
$ws.Cells.Item(2, 4).Value = [double]"0.9982866048812866"

$ws.Cells.Item(3, 1).Value = [double]"0.03866442292928696"
$ws.Cells.Item(3, 2).Value = [double]"0.9925417900085449"
$ws.Cells.Item(3, 3).Value = [double]"0.01010789256542921"
$ws.Cells.Item(3, 4).Value = [double]"0.9985314011573792"

$ws.Cells.Item(4, 1).Value = [double]"0.028123639523983"
$ws.Cells.Item(4, 2).Value = [double]"0.9924772977828979"
$ws.Cells.Item(4, 3).Value = [double]"0.004994395654648542"
$ws.Cells.Item(4, 4).Value = [double]"0.998653769493103"

$ws.Cells.Item(5, 1).Value = [double]"0.02112744934856892"
$ws.Cells.Item(5, 2).Value = [double]"0.9935089945793152"
$ws.Cells.Item(5, 3).Value = [double]"0.001030554762110114"
$ws.Cells.Item(5, 4).Value = [double]"0.9993880987167358"

$ws.Cells.Item(6, 1).Value = [double]"0.01390651520341635"
$ws.Cells.Item(6, 2).Value = [double]"0.996582567691803"
$ws.Cells.Item(6, 3).Value = [double]"0.001292318222112954"
$ws.Cells.Item(6, 4).Value = [double]"0.9992656707763672"

$ws.Cells.Item(7, 1).Value = [double]"0.01182360760867596"
$ws.Cells.Item(7, 2).Value = [double]"0.9968834519386292"
$ws.Cells.Item(7, 3).Value = [double]"0.00038805918302387"
$ws.Cells.Item(7, 4).Value = [double]"0.9998776316642761"

$ws.Cells.Item(8, 1).Value = [double]"0.01216852478682995"
$ws.Cells.Item(8, 2).Value = [double]"0.9970554113388062"
$ws.Cells.Item(8, 3).Value = [double]"0.0003844704478979111"
$ws.Cells.Item(8, 4).Value = [double]"0.9998776316642761"

$ws.Cells.Item(9, 1).Value = [double]"0.01036333851516247"
$ws.Cells.Item(9, 2).Value = [double]"0.9972488284111023"
$ws.Cells.Item(9, 3).Value = [double]"0.0001792005496099591"

$ws.Cells.Item(10, 1).Value = [double]"0.01013538055121899"
$ws.Cells.Item(10, 2).Value = [double]"0.9972058534622192"
$ws.Cells.Item(10, 3).Value = [double]"0.0002043738059001043"
$ws.Cells.Item(10, 4).Value = [double]"1"

$ws.Cells.Item(11, 1).Value = [double]"0.009414365515112877"
$ws.Cells.Item(11, 2).Value = [double]"0.9973993301391602"
$ws.Cells.Item(11, 3).Value = [double]"7.145854760892689E-05"

$ws.Cells.Item(12, 1).Value = [double]"0.01229505334049463"
$ws.Cells.Item(12, 2).Value = [double]"0.9964965581893921"
$ws.Cells.Item(12, 3).Value = [double]"7.885350350989029E-05"

$ws.Cells.Item(13, 1).Value = [double]"0.009517963975667953"
$ws.Cells.Item(13, 2).Value = [double]"0.9972273707389832"
$ws.Cells.Item(13, 3).Value = [double]"4.89922495034989E-05"
$ws.Cells.Item(13, 4).Value = [double]"1"

$ws.Cells.Item(14, 1).Value = [double]"0.009159478358924389"
$ws.Cells.Item(14, 2).Value = [double]"0.9975067973136902"
$ws.Cells.Item(14, 3).Value = [double]"7.038439071038738E-05"

$ws.Cells.Item(15, 1).Value = [double]"0.0102852089330554"
$ws.Cells.Item(15, 2).Value = [double]"0.996797502040863"
$ws.Cells.Item(15, 3).Value = [double]"2.73302666755626E-05"

$ws.Cells.Item(16, 1).Value = [double]"0.009235069155693054"
$ws.Cells.Item(16, 2).Value = [double]"0.9973778128623962"
$ws.Cells.Item(16, 3).Value = [double]"2.974271956190933E-05"

$ws.Cells.Item(17, 1).Value = [double]"0.009507066570222378"
$ws.Cells.Item(17, 2).Value = [double]"0.9971628785133362"
$ws.Cells.Item(17, 3).Value = [double]"2.342590960324742E-05"

$ws.Cells.Item(18, 1).Value = [double]"0.009460599161684513"
$ws.Cells.Item(18, 2).Value = [double]"0.9972273707389832"
$ws.Cells.Item(18, 3).Value = [double]"3.023721546924207E-05"

$ws.Cells.Item(19, 1).Value = [double]"0.01005534641444683"
$ws.Cells.Item(19, 2).Value = [double]"0.9970983862876892"
$ws.Cells.Item(19, 3).Value = [double]"2.633344229252543E-05"

$ws.Cells.Item(20, 1).Value = [double]"0.01048480812460184"
$ws.Cells.Item(20, 2).Value = [double]"0.9968189597129822"
$ws.Cells.Item(20, 3).Value = [double]"2.50649045483442E-05"

$ws.Cells.Item(21, 1).Value = [double]"0.009588534012436867"
$ws.Cells.Item(21, 2).Value = [double]"0.9971199035644531"
$ws.Cells.Item(21, 3).Value = [double]"2.577550003479701E-05"

$ws.Cells.Item(22, 1).Value = [double]"0.0083968136459589"
$ws.Cells.Item(22, 2).Value = [double]"0.9976786971092224"
$ws.Cells.Item(22, 3).Value = [double]"1.764850821928121E-05"

$ws.Cells.Item(23, 1).Value = [double]"0.009256658144295216"
$ws.Cells.Item(23, 2).Value = [double]"0.9972058534622192"
$ws.Cells.Item(23, 3).Value = [double]"5.560421413974836E-05"

$ws.Cells.Item(24, 1).Value = [double]"0.01049390714615583"
$ws.Cells.Item(24, 2).Value = [double]"0.99686199426651"
$ws.Cells.Item(24, 3).Value = [double]"6.46670232526958E-05"

$ws.Cells.Item(25, 1).Value = [double]"0.008817553520202637"
$ws.Cells.Item(25, 2).Value = [double]"0.9975712299346924"
$ws.Cells.Item(25, 3).Value = [double]"1.875059570011217E-05"

$ws.Cells.Item(26, 1).Value = [double]"0.009598582051694393"
$ws.Cells.Item(26, 2).Value = [double]"0.9971628785133362"
$ws.Cells.Item(26, 3).Value = [double]"1.629848884476814E-05"

$ws.Cells.Item(27, 1).Value = [double]"0.009154189378023148"
$ws.Cells.Item(27, 2).Value = [double]"0.9972273707389832"
$ws.Cells.Item(27, 3).Value = [double]"2.933229916379787E-05"

$ws.Cells.Item(28, 1).Value = [double]"0.0097627779468894"
$ws.Cells.Item(28, 2).Value = [double]"0.9970983862876892"
$ws.Cells.Item(28, 3).Value = [double]"3.89008127967827E-05"

$ws.Cells.Item(29, 1).Value = [double]"0.009429940022528172"
$ws.Cells.Item(29, 2).Value = [double]"0.9974637627601624"
$ws.Cells.Item(29, 3).Value = [double]"2.896509249694645E-05"

$ws.Cells.Item(30, 1).Value = [double]"0.01116597559303045"
$ws.Cells.Item(30, 2).Value = [double]"0.9967114925384521"
$ws.Cells.Item(30, 3).Value = [double]"6.803328869864345E-05"

$ws.Cells.Item(31, 1).Value = [double]"0.01012442074716091"
$ws.Cells.Item(31, 2).Value = [double]"0.9969264268875122"
$ws.Cells.Item(31, 3).Value = [double]"3.724498674273491E-05"

$ws.Cells.Item(32, 1).Value = [double]"0.008495848625898361"
$ws.Cells.Item(32, 2).Value = [double]"0.9975067973136902"
$ws.Cells.Item(32, 3).Value = [double]"2.836888779711444E-05"
$ws.Cells.Item(32, 4).Value = [double]"1"

$ws.Cells.Item(33, 1).Value = [double]"0.009771361947059631"
$ws.Cells.Item(33, 2).Value = [double]"0.9971843957901001"
$ws.Cells.Item(33, 3).Value = [double]"1.610372419236228E-05"

$ws.Cells.Item(34, 1).Value = [double]"0.009439610876142979"
$ws.Cells.Item(34, 2).Value = [double]"0.9971843957901001"
$ws.Cells.Item(34, 3).Value = [double]"2.006696013268083E-05"

$ws.Cells.Item(35, 1).Value = [double]"0.009268703870475292"
$ws.Cells.Item(35, 2).Value = [double]"0.9973778128623962"
$ws.Cells.Item(35, 3).Value = [double]"2.664950807229616E-05"

$ws.Cells.Item(36, 1).Value = [double]"0.009822673164308071"
$ws.Cells.Item(36, 2).Value = [double]"0.9970769286155701"
$ws.Cells.Item(36, 3).Value = [double]"4.070589420734905E-05"

$ws.Cells.Item(37, 1).Value = [double]"0.008154270239174366"
$ws.Cells.Item(37, 2).Value = [double]"0.9976786971092224"
$ws.Cells.Item(37, 3).Value = [double]"2.245294672320597E-05"

$ws.Cells.Item(38, 1).Value = [double]"0.009952053427696228"
$ws.Cells.Item(38, 2).Value = [double]"0.9972488284111023"
$ws.Cells.Item(38, 3).Value = [double]"2.867738658096641E-05"

$ws.Cells.Item(39, 1).Value = [double]"0.008208157494664192"
$ws.Cells.Item(39, 2).Value = [double]"0.9975927472114563"
$ws.Cells.Item(39, 3).Value = [double]"4.878155232290737E-05"

$ws.Cells.Item(40, 1).Value = [double]"0.008729130029678345"
$ws.Cells.Item(40, 2).Value = [double]"0.9974852800369263"
$ws.Cells.Item(40, 3).Value = [double]"3.926633144146763E-05"

$ws.Cells.Item(41, 1).Value = [double]"0.009076619520783424"
$ws.Cells.Item(41, 2).Value = [double]"0.9972488284111023"
$ws.Cells.Item(41, 3).Value = [double]"2.148236853827257E-05"

$ws.Cells.Item(42, 1).Value = [double]"0.008965450339019299"
$ws.Cells.Item(42, 2).Value = [double]"0.9972058534622192"
$ws.Cells.Item(42, 3).Value = [double]"2.146872975572478E-05"

$ws.Cells.Item(43, 1).Value = [double]"0.008881161920726299"
$ws.Cells.Item(43, 2).Value = [double]"0.9973778128623962"
$ws.Cells.Item(43, 3).Value = [double]"1.39682088047266E-05"

$ws.Cells.Item(44, 1).Value = [double]"0.00919934269040823"
$ws.Cells.Item(44, 2).Value = [double]"0.9970983862876892"
$ws.Cells.Item(44, 3).Value = [double]"1.024457469611662E-05"

$ws.Cells.Item(45, 1).Value = [double]"0.01101873908191919"
$ws.Cells.Item(45, 2).Value = [double]"0.9970983862876892"
$ws.Cells.Item(45, 3).Value = [double]"3.331856350996532E-05"

$ws.Cells.Item(46, 1).Value = [double]"0.009201523847877979"
$ws.Cells.Item(46, 2).Value = [double]"0.9974207878112793"
$ws.Cells.Item(46, 3).Value = [double]"4.427980093169026E-05"

$ws.Cells.Item(47, 1).Value = [double]"0.008552609011530876"
$ws.Cells.Item(47, 2).Value = [double]"0.9974423050880432"
$ws.Cells.Item(47, 3).Value = [double]"2.845989365596324E-05"

$ws.Cells.Item(48, 1).Value = [double]"0.008824082091450691"
$ws.Cells.Item(48, 2).Value = [double]"0.9972918629646301"
$ws.Cells.Item(48, 3).Value = [double]"2.785888864309527E-05"

$ws.Cells.Item(49, 1).Value = [double]"0.009276360273361206"
$ws.Cells.Item(49, 2).Value = [double]"0.9972058534622192"
$ws.Cells.Item(49, 3).Value = [double]"1.546728526591323E-05"

$ws.Cells.Item(50, 1).Value = [double]"0.009229958057403564"
$ws.Cells.Item(50, 2).Value = [double]"0.9971413612365723"
$ws.Cells.Item(50, 3).Value = [double]"9.153223800240085E-06"

$ws.Cells.Item(51, 1).Value = [double]"0.009109106846153736"
$ws.Cells.Item(51, 2).Value = [double]"0.9972703456878662"
$ws.Cells.Item(51, 3).Value = [double]"1.234460796695203E-05"
